# Backlog.xlsx - "Major code overhaul. Added menu system and revised
# component system." update:
#   - Row 11 (PhysicsComponent follow-up / sound hook-up item) and row 12
#     now have a "Date Completed" entry in column I.
#   - Row 11 grows taller (its Implementation Plan text wraps onto more
#     lines) so its row height becomes 90pt.
#   - The sheet's last saved selection moves to I13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Row 11 needs to be taller to fit its content (was the default 15pt).
$ws.Rows("11").RowHeight = 90

# Fill in the newly-tracked completion dates (column I) for rows 11 & 12.
$ws.Range("I11").Value = 40555
$ws.Range("I12").Value = 40555

# Leave the selection on I13, matching the saved sheet view.
$ws.Range("I13").Select() | Out-Null
